$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Matt Peter, Michael Riess, Jonah Kubath"
$ws.Range("B1").Select()
